$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("A5").Value = "TestCase_F4"
$ws.Range("B5").Value = "To verify that user receives a notification if someone likes his comment"
$ws.Range("C5").Value = "Y"
$ws.Range("D5").Value = "PASS"
